$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-09 Friday" "2024-02-10 Saturday"

Replace-Text "64×49=3136" "86×63=5418"
Replace-Text "25×97=2425" "53×69=3657"
Replace-Text "63×64=4032" "35×54=1890"
Replace-Text "87×17=1479" "24×26=624"
Replace-Text "52×15=780" "85×76=6460"

Replace-Text "26×80=2080" "99×23=2277"
Replace-Text "46×92=4232" "55×62=3410"
Replace-Text "68×55=3740" "48×81=3888"
Replace-Text "59×41=2419" "52×22=1144"
Replace-Text "37×58=2146" "96×36=3456"

Replace-Text "65×20=1300" "72×29=2088"
Replace-Text "60×97=5820" "88×33=2904"
Replace-Text "92×89=8188" "94×94=8836"
Replace-Text "17×26=442" "20×13=260"
Replace-Text "79×90=7110" "67×71=4757"

Replace-Text "82×85=6970" "38×71=2698"
Replace-Text "90×35=3150" "39×17=663"
Replace-Text "75×93=6975" "11×36=396"
Replace-Text "15×93=1395" "40×50=2000"
Replace-Text "14×85=1190" "95×27=2565"

Replace-Text "31×86=2666" "16×44=704"
Replace-Text "63×90=5670" "12×67=804"
Replace-Text "67×23=1541" "19×22=418"
Replace-Text "70×11=770" "77×14=1078"
Replace-Text "59×84=4956" "23×78=1794"
